# Updates cryptos list figures (price + volume(1h)) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while preventing Excel from
# auto-converting plain-numeric-looking strings (e.g. "578.71") into a
# numeric cell type. Non-numeric-looking strings (e.g. "66.712.95", which
# uses '.' as a thousands separator) are written as-is since Excel already
# keeps those as text.
function Set-TextValue($Cell, [string]$Text) {
    if ($Text -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

Set-TextValue $ws.Range("D2") "66.712.95"
$ws.Range("E2").Value = "  +2.39%  "

Set-TextValue $ws.Range("D3") "3.086.23"
$ws.Range("E3").Value = "  +4.58%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "578.71"
$ws.Range("E5").Value = "  +1.57%  "

Set-TextValue $ws.Range("D6") "168.05"
$ws.Range("E6").Value = "  +5.17%  "

$ws.Range("E7").Value = "  -0.07%  "

Set-TextValue $ws.Range("D8") "3.083.53"
$ws.Range("E8").Value = "  +4.76%  "

$ws.Range("E9").Value = "  +0.83%  "

Set-TextValue $ws.Range("D10") "6.58"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("E11").Value = "  +1.98%  "

$ws.Range("E12").Value = "  +5.41%  "

Set-TextValue $ws.Range("D13") "0.0000249"
$ws.Range("E13").Value = "  +1.92%  "

Set-TextValue $ws.Range("D14") "36.45"
$ws.Range("E14").Value = "  +6.99%  "

$ws.Range("E15").Value = "  -0.69%  "

Set-TextValue $ws.Range("D16") "3.595.85"
$ws.Range("E16").Value = "  +4.47%  "

Set-TextValue $ws.Range("D17") "66.733.70"
$ws.Range("E17").Value = "  +2.39%  "

Set-TextValue $ws.Range("D18") "7.22"
$ws.Range("E18").Value = "  +4.32%  "

Set-TextValue $ws.Range("D19") "3.084.83"
$ws.Range("E19").Value = "  +4.64%  "

Set-TextValue $ws.Range("D20") "16.21"
$ws.Range("E20").Value = "  +13.35%  "

Set-TextValue $ws.Range("D21") "466.74"
$ws.Range("E21").Value = "  +4.71%  "

Set-TextValue $ws.Range("D22") "0.715"
$ws.Range("E22").Value = "  +4.80%  "

Set-TextValue $ws.Range("D23") "7.55"
$ws.Range("E23").Value = "  +4.27%  "

Set-TextValue $ws.Range("D24") "83.18"
$ws.Range("E24").Value = "  +0.99%  "

Set-TextValue $ws.Range("D25") "2.33"
$ws.Range("E25").Value = "  +5.30%  "

Set-TextValue $ws.Range("D26") "12.86"
$ws.Range("E26").Value = "  +6.64%  "

Set-TextValue $ws.Range("D27") "10.15"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("E30").Value = "  +0.85%  "

Set-TextValue $ws.Range("D31") "2.67"
$ws.Range("E31").Value = "  +3.46%  "

$ws.Range("E32").Value = "  +0.57%  "

Set-TextValue $ws.Range("D33") "28.20"
$ws.Range("E33").Value = "  +3.99%  "

$ws.Range("E34").Value = "  +3.96%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  +2.71%  "

Set-TextValue $ws.Range("D37") "5.89"
$ws.Range("E37").Value = "  +3.65%  "

$ws.Range("E38").Value = "  +7.51%  "

Set-TextValue $ws.Range("D39") "47.02"
$ws.Range("E39").Value = "  +6.92%  "

# Rows 40/41 swap position: the coin previously in row 40 (OKB) moves to
# row 41, and the coin previously in row 41 (TheGraph) moves to row 40,
# each bringing its own refreshed price/volume figures.
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D40") "0.318"
$ws.Range("E40").Value = "  +6.59%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D41") "50.27"
$ws.Range("E41").Value = "  +2.70%  "

$ws.Range("E42").Value = "  +2.04%  "

Set-TextValue $ws.Range("D43") "8.69"
$ws.Range("E43").Value = "  +3.23%  "

Set-TextValue $ws.Range("D44") "2.83"
$ws.Range("E44").Value = "  -0.51%  "

Set-TextValue $ws.Range("D45") "0.0361"
$ws.Range("E45").Value = "  +2.85%  "

Set-TextValue $ws.Range("D46") "383.34"
$ws.Range("E46").Value = "  -0.60%  "

Set-TextValue $ws.Range("D47") "2.756.84"
$ws.Range("E47").Value = "  +1.39%  "

Set-TextValue $ws.Range("D48") "135.16"
$ws.Range("E48").Value = "  +1.73%  "

Set-TextValue $ws.Range("D50") "24.64"
$ws.Range("E50").Value = "  +6.25%  "

$ws.Range("E51").Value = "  +3.59%  "
